# Changing format and content of Conclusion in docx
#
# The "Conclusion" section heading paragraph gets an explicit direct
# character-formatting override: Times New Roman, 16pt (half-point value
# 32), applied both to the run that holds the "Conclusion" text and to the
# paragraph mark itself (so new text typed at the end of the heading
# inherits the same formatting).

$d = $word.ActiveDocument

# Locate the "Conclusion" heading. Use whole-word + case-sensitive match so
# this does not also hit the lowercase "conclusion" that appears later in
# the body text (e.g. "In conclusion, ...").
$rng = $d.Content
$found = $rng.Find.Execute("Conclusion", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $headingPara = $rng.Paragraphs(1)
    $headingRange = $headingPara.Range

    if ($headingRange.Text.TrimEnd([char]13, [char]7) -eq "Conclusion") {
        # Apply to the whole paragraph range (text run + paragraph mark) so
        # both the run's rPr and the paragraph mark's rPr (pPr/rPr) pick up
        # the new font/size, matching how Word formats a fully-selected
        # paragraph.
        $headingRange.Font.Name = "Times New Roman"
        $headingRange.Font.NameAscii = "Times New Roman"
        $headingRange.Font.NameOther = "Times New Roman"
        $headingRange.Font.NameBi = "Times New Roman"
        $headingRange.Font.Size = 16
        $headingRange.Font.SizeBi = 16
    }
}
